$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A36").Value = "Marco sartorelli"
$ws.Range("B36").Value = "Samuele Kettamier | SBARX"
$ws.Range("C36").Value = "Luca Barozzi | Modium"
$ws.Range("D36").Value = "Federico Andreis | iMontagna"
$ws.Range("E36").Value = "Maickol Azocar | MAI UNA GIOIA"
$ws.Range("F36").Value = "Francesco Vettori | F.C. Julia"
